$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -le 4; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("A$row")
    $oldValue = $cell.Text
    $newValue = $oldValue -replace "867746", "278857"
    $cell.Value = $newValue
}
